{"js": "// Add a new paragraph \"My name is Milad\" right after the existing\n// \"Hi\" paragraph (i.e. at the end of the document body), carrying\n// forward the same run/paragraph-mark formatting (en-GB language)\n// that Word would naturally inherit from the preceding paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"My name is Milad\", \"After\");\n\nawait context.sync();\n", "ps1": "# Add a new paragraph \"My name is Milad\" right after the existing\n# \"Hi\" paragraph (i.e. at the end of the document body), carrying\n# forward the same run/paragraph-mark formatting (en-GB language)\n# that Word would naturally inherit from the preceding paragraph.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"My name is Milad\"\n"}
